$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "admin"
$ws.Range("B2").Value = "manager"
$ws.Range("B1").Value = "password"

$ws.Range("D6").Select()
